$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table-number / year refresh ("tabel format baru": tables renumbered
# 4.2.3-4.2.5 -> 4.2.5-4.2.7, and the reporting year bumped 2020 -> 2021) ---

# H1 used to share its string with A1's "Tabel 4.2.3"; it now gets its own
# distinct "Tabel 4.2.5" label.
$ws.Range("H1").Value = "Tabel 4.2.5"

# P1 / W1 are rich-text ("Tabel" + " 4.2.x.") - rebuild the value then
# reapply the smaller run's own font (size 9, no underline) so the two
# runs stay distinct, matching the original " 4.2.x." run formatting.
$rngP1 = $ws.Range("P1")
$rngP1.Value = "Tabel 4.2.6."
$p1Run2 = $rngP1.Characters(6, 7)
$p1Run2.Font.Size = 9
$p1Run2.Font.Underline = $false

$rngW1 = $ws.Range("W1")
$rngW1.Value = "Tabel 4.2.7."
$w1Run2 = $rngW1.Characters(6, 7)
$w1Run2.Font.Size = 9
$w1Run2.Font.Underline = $false

# Title / subtitle cells: bump the survey year 2020 -> 2021.
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Wolo. 2021"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Wolo, 2021"
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Wolo, 2021"
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Wolo, 2021"

$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in Wolo Subdistrict, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village Wolo Subdistrict, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in Wolo Subdistrict, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in Wolo Subdistrict, 2021"

# Reset the view: scroll/selection back to the top-left cell (the saved
# file no longer pins the view at A10 with J21 selected).
$ws.Range("A1").Select()
